{"js": "// \"Added a few more slots\" - move the meta-description blurb from the top\n// of the document (right under the H1 title) down to the bottom, replacing\n// the AI image-generation prompt paragraph there, and add a bold restatement\n// of the page title just above it.\n\nconst body = context.document.body;\n\n// --- Step 1: remove the \"Meta description: ...\" paragraph that currently\n// sits directly under the title heading. ---\nconst paras = body.paragraphs;\nparas.load('items/text');\nawait context.sync();\n\nlet metaPara = null;\nfor (const p of paras.items) {\n  if (p.text.indexOf('Meta description') === 0) {\n    metaPara = p;\n    break;\n  }\n}\nif (metaPara) {\n  metaPara.delete();\n  await context.sync();\n}\n\n// --- Step 2: locate the final paragraph in the document (the one holding\n// the italic AI image-generation prompt) and insert a brand new bold\n// paragraph with the page title right before it. ---\nconst paras2 = body.paragraphs;\nparas2.load('items');\nawait context.sync();\nconst lastPara = paras2.items[paras2.items.length - 1];\n\n// Split off a fresh empty paragraph right before the last one...\nconst placeholder = lastPara.insertParagraph('', Word.InsertLocation.before);\nawait context.sync();\n\n// ...then stamp it with the exact run/formatting structure we want (a\n// leading empty run followed by a bold run) via a raw OOXML replace, so we\n// don't inherit the italic formatting of the neighboring paragraph.\nconst titleOoxml = '<?xml version=\"1.0\"?><pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\"><pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\"><pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Aztec Bonanza Free Slot Game | Review and Top Features</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\nplaceholder.insertOoxml(titleOoxml, 'Replace');\nawait context.sync();\n\n// --- Step 3: replace the text of the (now last) italic paragraph - the old\n// image-generation prompt - with the meta description copy, preserving its\n// run formatting (italic). ---\nconst searchText = \"Create a feature image for Aztec Bonanza that features a happy Maya warrior wearing glasses in a cartoon style. The image should be lively and exciting to match the adventurous theme of the game. You can incorporate elements like gemstones, the totem pole, or the temple in the background to give the image an Aztec touch. Make sure to highlight the cascading game mechanics and the free spins and giant symbol features to entice players. The overall image should convey a feeling of excitement and adventure to entice players to try this game.\";\nconst results = body.search(searchText, { matchCase: true, matchWholeWord: false });\nresults.load('items');\nawait context.sync();\nif (results.items.length > 0) {\n  results.items[0].insertText('Discover the exciting gameplay of Aztec Bonanza for free. Learn about its features, pros, and cons. Play now!', Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# \"Added a few more slots\" - move the meta-description blurb from the top\n# of the document (right under the H1 title) down to the bottom, replacing\n# the AI image-generation prompt paragraph there, and add a bold restatement\n# of the page title just above it.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: remove the \"Meta description: ...\" paragraph that currently\n# sits directly under the title heading. ---\n$metaPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"Meta description*\") {\n        $metaPara = $p\n        break\n    }\n}\nif ($metaPara -ne $null) {\n    $metaPara.Range.Delete()\n}\n\n# --- Step 2: locate the final paragraph in the document (the one holding\n# the italic AI image-generation prompt) and insert a brand new bold\n# paragraph with the page title right before it. ---\n$lastIdx = $d.Paragraphs.Count\n$imgPara = $d.Paragraphs($lastIdx)\n$imgPara.Range.InsertParagraphBefore()\n\n$titleIdx = $lastIdx\n$titlePara = $d.Paragraphs($titleIdx)\n$titleXml = '<w:p xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Aztec Bonanza Free Slot Game | Review and Top Features</w:t></w:r></w:p>'\n$titlePara.Range.InsertXML($titleXml)\n\n# --- Step 3: replace the text of the (now last) italic paragraph - the old\n# image-generation prompt - with the meta description copy. ---\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Execute(\"Create a feature image for Aztec Bonanza that features a happy Maya warrior wearing glasses in a cartoon style. The image should be lively and exciting to match the adventurous theme of the game. You can incorporate elements like gemstones, the totem pole, or the temple in the background to give the image an Aztec touch. Make sure to highlight the cascading game mechanics and the free spins and giant symbol features to entice players. The overall image should convey a feeling of excitement and adventure to entice players to try this game.\", $false, $false, $false, $false, $false, $true, 1, $false, \"Discover the exciting gameplay of Aztec Bonanza for free. Learn about its features, pros, and cons. Play now!\", 2)\n"}
